# =====================================================================
# Atualização do sistema de registros com melhorias em gráficos e orçamentos
# =====================================================================
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheet "Registros": clear stale data rows and write a single new
#    record row (row 2). Leading-zero / numeric-looking text values
#    (Data, ID, Número Projeto) are forced to Text so Excel keeps them
#    as literal strings instead of silently converting them to numbers
#    or dates; the style is reset back to "Normal" afterwards so no
#    extra formatting is left behind on the cell.
# ---------------------------------------------------------------------
$wsReg = $wb.Worksheets.Item("Registros")
$wsReg.Range("A2:I16").Clear()

$wsReg.Range("A2").NumberFormat = "@"
$wsReg.Range("A2").Value = "2025-10-24"
$wsReg.Range("A2").Style = "Normal"

$wsReg.Range("B2").NumberFormat = "@"
$wsReg.Range("B2").Value = "010"
$wsReg.Range("B2").Style = "Normal"

$wsReg.Range("C2").Value = "JOSE GENILSOS MARTINS SOARES"
$wsReg.Range("D2").Value = "Marcenaria Estrutural"
$wsReg.Range("E2").Value = "ODS"

$wsReg.Range("F2").NumberFormat = "@"
$wsReg.Range("F2").Value = "12"
$wsReg.Range("F2").Style = "Normal"

$wsReg.Range("G2").Value = "15:50"
$wsReg.Range("H2").Value = "16:05"
$wsReg.Range("I2").Value = "registro"

# ---------------------------------------------------------------------
# 2) Sheet "Gráficos": widen the columns, apply the header style (reuse
#    the same bold+centered style already used by "Registros" header so
#    no superfluous style entries are created), and fill in the budget
#    summary rows.
# ---------------------------------------------------------------------
$wsGraf = $wb.Worksheets.Item("Gráficos")

$wsGraf.Columns.Item(1).ColumnWidth = 29.1666666666667
$wsGraf.Columns.Item(2).ColumnWidth = 14.1666666666667
$wsGraf.Columns.Item(3).ColumnWidth = 14.1666666666667
$wsGraf.Columns.Item(4).ColumnWidth = 14.1666666666667
$wsGraf.Columns.Item(5).ColumnWidth = 14.1666666666667

$wsReg.Range("A1").Copy()
$wsGraf.Range("A1:E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$grafData = @(
  @("Marcenaria Estrutural - CDS - 12", 0, 30, 30),
  @("Marcenaria Estrutural - CDS - 15", 0, 65, 65),
  @("Marcenaria Estrutural - FÁBRICA - 12", 0, 23, 23),
  @("Marcenaria Estrutural - ODS - 12", 0, 23, 23),
  @("Marcenaria Estrutural - TDS - 15", 0, 65, 65),
  @("Marcenaria Móvel - CDS - 12", 0, 23, 23),
  @("Marcenaria Móvel - ODS - 12", 0, 30, 30)
)

$r = 2
foreach ($row in $grafData) {
  $wsGraf.Range("A$r").Value = $row[0]
  $wsGraf.Range("B$r").NumberFormat = "0.00"
  $wsGraf.Range("B$r").Value = $row[1]
  $wsGraf.Range("C$r").NumberFormat = "0"
  $wsGraf.Range("C$r").Value = $row[2]
  $wsGraf.Range("D$r").NumberFormat = "0.00"
  $wsGraf.Range("D$r").Value = $row[3]
  $r++
}

# ---------------------------------------------------------------------
# 3) New sheet "Orçamentos" at the end of the workbook, listing the
#    budgeted hours per Área/Projeto/Número Projeto.
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsOrc = $wb.Worksheets.Add($null, $lastSheet)
$wsOrc.Name = "Orçamentos"

$wsOrc.Columns.Item(1).ColumnWidth = 19.1666666666667
$wsOrc.Columns.Item(2).ColumnWidth = 19.1666666666667
$wsOrc.Columns.Item(3).ColumnWidth = 19.1666666666667
$wsOrc.Columns.Item(4).ColumnWidth = 14.1666666666667

$wsReg.Range("A1").Copy()
$wsOrc.Range("A1:D1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wsOrc.Range("A1").Value = "Área"
$wsOrc.Range("B1").Value = "Projeto"
$wsOrc.Range("C1").Value = "Número Projeto"
$wsOrc.Range("D1").Value = "Horas Orçadas"

$orcData = @(
  @("Marcenaria Estrutural", "CDS", "15", 65),
  @("Marcenaria Móvel", "ODS", "12", 30),
  @("Marcenaria Estrutural", "CDS", "12", 30),
  @("Marcenaria Estrutural", "TDS", "15", 65),
  @("Marcenaria Estrutural", "ODS", "12", 23),
  @("Marcenaria Estrutural", "FÁBRICA", "12", 23),
  @("Marcenaria Móvel", "CDS", "12", 23)
)

$r = 2
foreach ($row in $orcData) {
  $wsOrc.Range("A$r").Value = $row[0]
  $wsOrc.Range("B$r").Value = $row[1]

  $wsOrc.Range("C$r").NumberFormat = "@"
  $wsOrc.Range("C$r").Value = $row[2]
  $wsOrc.Range("C$r").Style = "Normal"

  $wsOrc.Range("D$r").NumberFormat = "0"
  $wsOrc.Range("D$r").Value = $row[3]
  $r++
}

# ---------------------------------------------------------------------
# Keep "Registros" as the active sheet/tab, as in the original workbook.
# ---------------------------------------------------------------------
$wsReg.Activate() | Out-Null
$wsReg.Range("A1").Select() | Out-Null
